# LOQ4252.docx edit
#
# Net effect (per the target diff): several text blocks rotate between
# paragraphs, the old "Programa resumido" body + the "Programa" heading
# that used to sit right under it are removed from there and re-created
# (heading + a new body line) right after the big numbered "Programa"
# body, the Método/Critério/Norma bullet's three content runs are
# rewritten, and the Bibliografia body is replaced by the docente line.
#
# Structural changes (paragraph delete/insert) are done first, while the
# anchor text used to find them is still unique in the document; the
# plain text substitutions happen afterwards.

$d = $word.ActiveDocument
$VT = [char]11

# ---------------------------------------------------------------------
# Step 1: remove the old "Conceitos ligados..." paragraph and the
# "Programa" heading that directly followed "Programa resumido".
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Programa resumido`r") {
        $toDelete = $d.Range($d.Paragraphs.Item($i + 1).Range.Start, $d.Paragraphs.Item($i + 2).Range.End)
        $toDelete.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# Step 2: right after the detailed numbered "Programa" body, insert a
# new "Programa" Heading 2 paragraph followed by a new body paragraph
# holding the old Método sentence.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.ParagraphStyle.NameLocal -eq "Normal" -and $para.Range.Text -like "1. Conceitos ligados ao escoamento de fluídos*") {
        $para.Range.InsertParagraphAfter()
        $headPara = $d.Paragraphs.Item($i + 1)
        $headPara.Range.Text = "Programa"
        $headPara.Range.Style = "Heading 2"

        $headPara.Range.InsertParagraphAfter()
        $bodyPara = $d.Paragraphs.Item($i + 2)
        $bodyPara.Range.Style = "Normal"
        $bodyPara.Range.Text = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios, aulas de laboratório."
        break
    }
}

# ---------------------------------------------------------------------
# Step 3: Objetivos body <- old "Programa resumido" body text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Fornecer os conceitos básicos de Mecânica dos Fluidos e Transferência de Calor e Massa com aplicações à Engenharia. Possibilitar aos alunos uma base científica para que possam se desenvolver em demais disciplinas tecnológicas do curso.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Conceitos ligados ao escoamento de fluídos e equações fundamentais, Escoamento incompressível de fluidos não viscosos, Escoamento viscoso incompressível, Transferência de Calor. Transferência de Massa",
    2) | Out-Null

# ---------------------------------------------------------------------
# Step 4: Docente(s) ListBullet body <- old Objetivos body text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "4808662 - Lucrécio Fábio dos Santos",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Fornecer os conceitos básicos de Mecânica dos Fluidos e Transferência de Calor e Massa com aplicações à Engenharia. Possibilitar aos alunos uma base científica para que possam se desenvolver em demais disciplinas tecnológicas do curso.",
    2) | Out-Null

# ---------------------------------------------------------------------
# Step 5: rewrite the Avaliação bullet's three content runs (the text
# after each bold label), working back-to-front so earlier Range
# offsets captured via Find stay valid.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.ParagraphStyle.NameLocal -eq "List Bullet" -and $para.Range.Text -like "Método: *Norma de recuperação: *") {
        $mLabel = $d.Range($para.Range.Start, $para.Range.End)
        $mLabel.Find.Execute("Método: ") | Out-Null

        $cLabel = $d.Range($para.Range.Start, $para.Range.End)
        $cLabel.Find.Execute("Critério: ") | Out-Null

        $nLabel = $d.Range($para.Range.Start, $para.Range.End)
        $nLabel.Find.Execute("Norma de recuperação: ") | Out-Null

        $content3 = $d.Range($nLabel.End, $para.Range.End)
        $content3.Text = "1. FOX, R.W., MCDONALD, A.T., “Introdução à Mecânica dos Fluidos”, Ed. Guanabara Koogan." + $VT + `
            "2. STREETER, V.L., WYLE,E.B., “Mecânica dos Fluidos”, Ed. Mc Graw Hill." + $VT + `
            "3. OZISIK,M.N., “Transferência de Calor.”, Ed. Guanabara Koogan." + $VT + `
            "4. INCROPERA, F.P.W., “Fundamentos de Transferência de Calor e Massa”, Ed. Guanabara Koogan." + $VT + `
            "5. MUNSON, B.R.; YOUNG, D.F.; OKIISHI, T.H. Fundamentos da Mecânica dos Fluidos. Editora Edgard Blucher" + $VT + `
            "6 - GIORGETI, M. (2012) Fundamentos de Fenômenos de Transporte. Editora Campus"

        $content2 = $d.Range($cLabel.End, $nLabel.Start)
        $content2.Text = "Aplicação de uma prova envolvendo o assunto de todo semestre." + $VT + "NR (nota da recuperação) = (M1 + NR)/2." + $VT

        $content1 = $d.Range($mLabel.End, $cLabel.Start)
        $content1.Text = "Nota de duas provas (P1 e P2)" + $VT + "Fórmula: M1 = (P1 + 2 x P2)/3.." + $VT
        break
    }
}

# ---------------------------------------------------------------------
# Step 6: Bibliografia body <- old Docente(s) ListBullet text (whole
# paragraph content replaced).
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.ParagraphStyle.NameLocal -eq "Normal" -and $para.Range.Text -like "1. FOX, R.W., MCDONALD*") {
        $para.Range.Text = "4808662 - Lucrécio Fábio dos Santos"
        break
    }
}
